$wb = $excel.ActiveWorkbook

# 展览 (Exhibition) sheet updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1367
$ws1.Range("G2").Value = 69
$ws1.Range("F3").Value = 1443
$ws1.Range("F7").Value = 679
$ws1.Range("F8").Value = 119
$ws1.Range("F11").Value = 2476
$ws1.Range("F13").Value = 1507
$ws1.Range("F14").Value = 313
$ws1.Range("G14").Value = "已售罄"
$ws1.Range("F15").Value = 248
$ws1.Range("F16").Value = 620
$ws1.Range("F17").Value = 794
$ws1.Range("F18").Value = 84
$ws1.Range("F22").Value = 30
$ws1.Range("F24").Value = 5088
$ws1.Range("F26").Value = 553
$ws1.Range("F27").Value = 83
$ws1.Range("F28").Value = 161
$ws1.Range("F31").Value = 225
$ws1.Range("F32").Value = 33
$ws1.Range("F33").Value = 1043
$ws1.Range("F34").Value = 743
$ws1.Range("F36").Value = 54
$ws1.Range("F39").Value = 1076
$ws1.Range("F42").Value = 175
$ws1.Range("F44").Value = 54

# 演出 (Performance) sheet updates
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F7").Value = 9

# 全部类型 (All Types) sheet updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1367
$ws4.Range("G2").Value = 69
$ws4.Range("F5").Value = 1443
$ws4.Range("F11").Value = 679
$ws4.Range("F12").Value = 119
$ws4.Range("F14").Value = 9
$ws4.Range("F17").Value = 2476
$ws4.Range("F19").Value = 1507
$ws4.Range("F20").Value = 313
$ws4.Range("G20").Value = "已售罄"
$ws4.Range("F21").Value = 248
$ws4.Range("F22").Value = 620
$ws4.Range("F24").Value = 794
$ws4.Range("F25").Value = 84
$ws4.Range("F28").Value = 30
$ws4.Range("F29").Value = 5088
$ws4.Range("F31").Value = 553
$ws4.Range("F32").Value = 83
$ws4.Range("F33").Value = 161
$ws4.Range("F36").Value = 225
$ws4.Range("F37").Value = 33
$ws4.Range("F38").Value = 1043
$ws4.Range("F39").Value = 743
$ws4.Range("F40").Value = 54
$ws4.Range("F42").Value = 1076
$ws4.Range("F44").Value = 175
$ws4.Range("F46").Value = 54

